$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add I1 = "I0", J1 = "IF" with the same style as the
# --- existing bold/bordered header cells (copy format from H1, then set text).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-25: add new I/J numeric values.
$data = @(
    @{ Row = 2;  I = 8;  J = 9 },
    @{ Row = 3;  I = 6;  J = 6 },
    @{ Row = 4;  I = 7;  J = 8 },
    @{ Row = 5;  I = 5;  J = 7 },
    @{ Row = 6;  I = 1;  J = 3 },
    @{ Row = 7;  I = 7;  J = 7 },
    @{ Row = 8;  I = 8;  J = 9 },
    @{ Row = 9;  I = 10; J = 10 },
    @{ Row = 10; I = 1;  J = 3 },
    @{ Row = 11; I = 7;  J = 8 },
    @{ Row = 12; I = 7;  J = 8 },
    @{ Row = 13; I = 5;  J = 8 },
    @{ Row = 14; I = 7;  J = 9 },
    @{ Row = 15; I = 11; J = 11 },
    @{ Row = 16; I = 8;  J = 8 },
    @{ Row = 17; I = 5;  J = 6 },
    @{ Row = 18; I = 8;  J = 8 },
    @{ Row = 19; I = 9;  J = 9 },
    @{ Row = 20; I = 6;  J = 6 },
    @{ Row = 21; I = 6;  J = 7 },
    @{ Row = 22; I = 3;  J = 4 },
    @{ Row = 23; I = 4;  J = 5 },
    @{ Row = 24; I = 6;  J = 6 },
    @{ Row = 25; I = 4;  J = 4 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I
    $ws.Cells.Item($item.Row, 10).Value = $item.J
}
